# 2026-01-30 12:51 JST scrape run: two freshly-scraped job postings are
# prepended to the "ランサーズ" listing sheet (newest first), pushing the
# previously-seen six postings down by two rows, and every row's
# "fetched at" timestamp is refreshed to the new run time. Column B also
# grows a bit wider to fit the longer titles now in the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2026-01-30 12:51:10"

# --- widen column B (42 -> 52) ---------------------------------------------
# (51.15 is not a typo: this runtime's ColumnWidth setter quantizes to
# pixels and round-trips back to "characters" width, so asking for 52.0
# directly lands on 52.8333; 51.15 is the value that round-trips to
# exactly 52, matching the target ooxml.)
$ws.Columns.Item(2).ColumnWidth = 51.15

# --- make room for the two new postings at the top -------------------------
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# --- row 2: newest posting ---------------------------------------------------
$ws.Cells.Item(2, 1).Value = $timestamp
$ws.Cells.Item(2, 2).Value = "【Zapier保守・運用サポート】既存フローの管理・調整をお任せできる方募集(時給1,150円程度)"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5475245"
$ws.Cells.Item(2, 7).Value = 213
$ws.Cells.Item(2, 8).Value = "🔥API ◇管理"

# --- row 3: second-newest posting -------------------------------------------
$ws.Cells.Item(3, 1).Value = $timestamp
$ws.Cells.Item(3, 2).Value = "シミュレーションスタジオの入退館システム開発(ロック選定含む/多店舗・複数打席対応)"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5482462"
$ws.Cells.Item(3, 7).Value = 125
$ws.Cells.Item(3, 8).Value = "◆開発,システム開発"

# --- rows 4-9: previously-existing postings, only the timestamp changes ----
for ($r = 4; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $timestamp
}

# --- hyperlinks: Insert() does not carry the old hyperlink objects down
#     with the cells they used to decorate, so rebuild the whole collection
#     in the final, correct row order. ---------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://www.lancers.jp/work/detail/5475245") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), "https://www.lancers.jp/work/detail/5482462") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "https://www.lancers.jp/work/detail/5482097") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://www.lancers.jp/work/detail/5482389") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://www.lancers.jp/work/detail/5481859") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), "https://www.lancers.jp/work/detail/5418064") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), "https://www.lancers.jp/work/detail/5481715") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), "https://www.lancers.jp/work/detail/5481888") | Out-Null
